$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the absolute path recorded for the workbook (x15ac:absPath), best effort
$wb.Path = "C:\SKL\excel\"
$wb.FullName = "C:\SKL\excel\L9731_底稿_人工檢核表1.xlsx"

# Give the new column its own default formatting/width like the other
# right-aligned data columns (A, D, E, F).
$ws.Columns.Item(9).HorizontalAlignment = -4152  # xlRight
$ws.Columns.Item(9).VerticalAlignment = -4108    # xlCenter
$ws.Columns.Item(9).ColumnWidth = 14.285714285714286  # renders as width 15 in OOXML

# Add new header "不足額金額" in column I, matching the formatting of the
# neighboring header cell (H1) - same font/fill/border/alignment - applied
# last so the header cell keeps the shared header style (not the column
# default).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "不足額金額"

# Update the active selection to I2 (matches the recorded cursor position)
$ws.Range("I2").Select()
